# Generate Report for Handoff
# Adds a new localization-status row (for file
# f03db06f-0e4b-497e-aaad-79a5b43bce82.md) to the Overview, zh-cn and
# de-de worksheets, matching the layout of the existing 62f5f371... row.

$wb = $excel.ActiveWorkbook

$newFileBase   = "f03db06f-0e4b-497e-aaad-79a5b43bce82"
$newFileName   = "$newFileBase.md"
$newFileDisp   = "e2e\$newFileBase.md"
$hyperlinkBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/670306ea8673fde87a84028ae3947555f8bbc197/e2e/$newFileName"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "'$newFileName"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkBase, [System.Type]::Missing, [System.Type]::Missing, $newFileDisp)
$wsOverview.Range("C3").Value = "'.md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = "'Ready for handoff"
$wsOverview.Range("F3").Value = "'Ready for handoff"
$wsOverview.Range("G3").Value = "'2016-08-12 12:45:25"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhXlf  = "$newFileBase.719faabd246b69d8128b49e5c96a548c2cbb093b.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkBase, [System.Type]::Missing, [System.Type]::Missing, $newFileName)
$wsZhCn.Range("B3").Value = "'.md"
$wsZhCn.Range("C3").Value = "'Ready for handoff"
$wsZhCn.Range("D3").Value = "'e2e"
$wsZhCn.Range("E3").Value = "'ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "'$zhXlf"
$wsZhCn.Range("H3").Value = "'2016-08-12 12:45:18"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I3").Value = "'"
$wsZhCn.Range("J3").Value = "'"
$wsZhCn.Range("K3").Value = "'0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$deXlf  = "$newFileBase.719faabd246b69d8128b49e5c96a548c2cbb093b.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkBase, [System.Type]::Missing, [System.Type]::Missing, $newFileName)
$wsDeDe.Range("B3").Value = "'.md"
$wsDeDe.Range("C3").Value = "'Ready for handoff"
$wsDeDe.Range("D3").Value = "'e2e"
$wsDeDe.Range("E3").Value = "'ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "'$deXlf"
$wsDeDe.Range("H3").Value = "'2016-08-12 12:45:25"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I3").Value = "'"
$wsDeDe.Range("J3").Value = "'"
$wsDeDe.Range("K3").Value = "'0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

Write-Host "Added handoff row for $newFileName to Overview, zh-cn and de-de sheets."
